# Weekly data refresh for the Hortaliza / Vega Central Mapocho de Santiago - Sandia sheet.
# The whole data table (rows 114..196) gets shifted down by one row (a new day's
# reading is inserted at the top of that block, row 114), and the row that falls
# off the bottom (old row 196) becomes a brand-new row 197.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture the last existing data row (196) in full BEFORE we shift anything,
#    since that row's original contents become the new row 197.
$lastRow = $ws.Range("A196:R196").Value2

# 2) Shift the variable part of the table (rows 114..195) down by one row, into
#    115..196. Columns A,B,C,E,F,G,H,Q,R are constant across the whole table, so
#    only D (Fecha) and I:P (Calidad..Precio $/Kg) need to move.
$datesBlock = $ws.Range("D114:D195").Value2
$ws.Range("D115:D196").Value2 = $datesBlock

$dataBlock = $ws.Range("I114:P195").Value2
$ws.Range("I115:P196").Value2 = $dataBlock

# 3) Append the captured old row 196 as the new row 197.
$ws.Range("A197:R197").Value2 = $lastRow
# Row 197 is brand new, so it doesn't inherit the date number format that the
# "Fecha" column (D) uses elsewhere in the table - restore it explicitly.
$ws.Range("D197").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 4) Overwrite row 114 with the newly-reported reading (a new "Perú" origin entry).
$ws.Range("D114").Value2 = 44455
$ws.Range("J114").Value2 = 160
$ws.Range("K114").Value2 = 1000
$ws.Range("L114").Value2 = 1200
$ws.Range("M114").Value2 = 1100
$ws.Range("O114").Value2 = "Perú"
$ws.Range("P114").Value2 = 1100
